$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - set values then copy formatting (bold font + border + alignment) from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for columns I and J
$data = @{
    2  = @(10, 10)
    3  = @(7, 8)
    4  = @(9, 9)
    5  = @(6, 7)
    6  = @(9, 9)
    7  = @(6, 6)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(5, 5)
    11 = @(10, 10)
    12 = @(5, 6)
    13 = @(8, 8)
    14 = @(2, 2)
    15 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
